$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 113 (this shifts the existing rows 113-150 down to 114-151,
# preserving all of their data/formatting - matching the "old row N becomes new row N+1" pattern
# seen across the whole diff).
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,R are identical to the rest of the Perejil/Feria Lagunitas block.
$ws.Cells.Item(113, 1).Value = 4
$ws.Cells.Item(113, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(113, 3).Value = "Los Lagos"
$ws.Cells.Item(113, 4).Value = 44463
$ws.Cells.Item(113, 5).Value = 10
$ws.Cells.Item(113, 6).Value = 100112044
$ws.Cells.Item(113, 7).Value = "Perejil"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 180
$ws.Cells.Item(113, 11).Value = 4500
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = 4500
$ws.Cells.Item(113, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(113, 15).Value = "Región Metropolitana"
$ws.Cells.Item(113, 16).Value = 1500
$ws.Cells.Item(113, 17).Value = 3
$ws.Cells.Item(113, 18).Value = "Hortaliza"
